$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the monthly cochera (garage) value
$ws.Range("A1").Value = 1300

# Recalculate dependent formulas (A2 = A1*12, A3 = B2/A2-1)
$excel.Calculate()

# Move the active selection to A2, matching the saved view state
$ws.Range("A2").Select()
